$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.460.95'
$ws.Range("E2").Value = '  -2.51%  '
$ws.Range("D3").Value = '2.893.77'
$ws.Range("E3").Value = '  -3.71%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.43'
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.38'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.503'
$ws.Range("E8").Value = '  -2.66%  '
$ws.Range("D9").Value = '2.889.95'
$ws.Range("E9").Value = '  -3.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.62'
$ws.Range("E10").Value = '  +6.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.144'
$ws.Range("E11").Value = '  -3.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.445'
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("E13").Value = '  -3.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.16'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '3.374.40'
$ws.Range("E16").Value = '  -3.59%  '
$ws.Range("D17").Value = '60.438.29'
$ws.Range("E17").Value = '  -2.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.79'
$ws.Range("E18").Value = '  -2.76%  '
$ws.Range("D19").Value = '2.896.42'
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '425.41'
$ws.Range("E20").Value = '  -4.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.58'
$ws.Range("E21").Value = '  -3.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.668'
$ws.Range("E22").Value = '  -2.57%  '
$ws.Range("E23").Value = '  -3.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.66'
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.01'
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.16'
$ws.Range("E26").Value = '  -3.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.77'
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.21'
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.61'
$ws.Range("E32").Value = '  -3.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.51'
$ws.Range("E33").Value = '  -3.26%  '
$ws.Range("E34").Value = '  -3.96%  '
$ws.Range("D35").Value = '0.0₃0836'
$ws.Range("E35").Value = '  -1.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -1.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.65'
$ws.Range("E37").Value = '  -2.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.68'
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.123'
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.15'
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '371.67'
$ws.Range("E45").Value = '  -5.50%  '
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '133.38'
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("D48").Value = '2.646.80'
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.09'
$ws.Range("E50").Value = '  +5.84%  '
$ws.Range("E51").Value = '  -0.88%  '
